$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 4 is "Working Demo : <url>" where the url is a separate
# hyperlinked run. Replace just the url portion, preserving the run's
# own formatting (hyperlink, underline, color, fonts).
$label = "Working Demo : "
$para4 = $tr.Paragraphs(4, 1)
$urlStart = $para4.Start + $label.Length
$urlLen = $para4.Length - $label.Length
$urlRange = $tr.Characters($urlStart, $urlLen)
$urlRange.Text = "https://www.youtube.com/watch?v=S-oK19eegDg"

# Remove the now-redundant empty paragraph that used to follow the
# "Working Demo" line (paragraph index 5).
$blankPara = $tr.Paragraphs(5, 1)
$blankPara.Delete()
